$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------------
# "BLK" sheet was updated with a newly-run "Wire and facade" experiment result
# (inserted ahead of the previous "no noise" stats block) and a newly-run
# "Veg and Facade" experiment result (also inserted ahead of its previous
# "no noise" stats block). The older stats blocks are pushed further down the
# sheet, and the "Wire and facade" section header moves from row 13 to row 21.
# ----------------------------------------------------------------------------

# --- Push the existing "no noise" stats for Veg and Facade down (rows 4-11 -> 12-17) ---
$ws.Range('C12').Value2 = 'The number of total instances in the two classes: 11819'
$ws.Range('C13').Value2 = 'The number of correctlly classified: 10955'
$ws.Range('C14').Value2 = 'The accuracy of classification: 92.69 %'
$ws.Range('C15').Value2 = 'The confusion matrix:'
$ws.Range('C16').Value2 = '[[6397  344]'
$ws.Range('C17').Value2 = ' [ 520 4558]]'

# --- Push the existing "Wire and facade" block down (rows 13-23 -> 21-37) ---
$ws.Range('B21').Value2 = 'Wire and facade'
$ws.Range('C22').Value2 = 'no noise'
$ws.Range('F22').Value2 = 'noise0.5'

$ws.Range('C32').Value2 = 'The number of total instances in the two classes: 6053'
$ws.Range('C33').Value2 = 'The number of correctlly classified: 5611'
$ws.Range('C34').Value2 = 'The accuracy of classification: 92.70 %'
$ws.Range('C35').Value2 = 'The confusion matrix:'
$ws.Range('C36').Value2 = '[[ 709  266]'
$ws.Range('C37').Value2 = ' [ 176 4902]]'

# --- Insert the new "Wire and facade" result block (rows 24-31) ---
$ws.Range('C24').Value2 = 'Training time without validation: 0.328284 '
$ws.Range('F24').Value2 = 'Training time without validation: 0.290879 '
$ws.Range('C25').Value2 = 'The number of total instances in the two classes: 13948'
$ws.Range('F25').Value2 = 'Test time for all test data: 0.282190 '
$ws.Range('C26').Value2 = 'The number of correctlly classified: 13003'
$ws.Range('F26').Value2 = 'The number of total instances in the two classes: 6079'
$ws.Range('C27').Value2 = 'The accuracy of classification: 93.22 %'
$ws.Range('F27').Value2 = 'The number of correctlly classified: 5642'
$ws.Range('C28').Value2 = 'The confusion matrix:'
$ws.Range('F28').Value2 = 'The accuracy of classification: 92.81 %'
$ws.Range('C29').Value2 = '[[ 1607   625]'
$ws.Range('F29').Value2 = 'The confusion matrix:'
$ws.Range('C30').Value2 = ' [  320 11396]]'
$ws.Range('F30').Value2 = '[[ 719  265]'
$ws.Range('C31').Value2 = 'Test time for all test data: 0.302620 '
$ws.Range('F31').Value2 = ' [ 172 4923]]'

# --- Insert the new "Veg and Facade" result block (rows 4-11) ---
$ws.Range('C4').Value2 = 'Training time without validation: 0.444265 '
$ws.Range('C5').Value2 = 'The number of total instances in the two classes: 27432'
$ws.Range('C6').Value2 = 'The number of correctlly classified: 25462'
$ws.Range('C7').Value2 = 'The accuracy of classification: 92.82 %'
$ws.Range('C8').Value2 = 'The confusion matrix:'
$ws.Range('C9').Value2 = '[[14880   836]'
$ws.Range('C10').Value2 = ' [ 1134 10582]]'
$ws.Range('C11').Value2 = 'Test time for all test data: 0.256323 '

# --- Clear the cells left behind by the row shifts ---
$ws.Range('B13').ClearContents()
$ws.Range('F14').ClearContents()
$ws.Range('F16').ClearContents()
$ws.Range('F17').ClearContents()
$ws.Range('C18').ClearContents()
$ws.Range('F18').ClearContents()
$ws.Range('C19').ClearContents()
$ws.Range('F19').ClearContents()
$ws.Range('C20').ClearContents()
$ws.Range('F20').ClearContents()
$ws.Range('C21').ClearContents()
$ws.Range('F21').ClearContents()
$ws.Range('C23').ClearContents()
$ws.Range('F23').ClearContents()

# --- Match final cursor / window state from the author's session ---
$w = $wb.Windows.Item(1)
$w.Left = 3260
$w.Top = 1120

$ws.Range('C16').Select()
